$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Formula")

# Row 2: new cells D2, I2
$ws.Range("D2").Formula = "=30*12"
$ws.Range("I2").Formula = "=PI()"

# Row 3: new cell D3
$ws.Range("D3").Formula = "=365/2"

# Row 4: B4 formula changes; add G4, H4 labels
$ws.Range("B4").Formula = "=PI()"
$ws.Range("G4").Value = "day length"
$ws.Range("H4").Value = "night length"

# Row 5 (new row of labels + MAX formulas)
$ws.Range("F5").Value = "max"
$ws.Range("G5").Formula = "=MAX(C$7:C$371)"
$ws.Range("H5").Formula = "=MAX(D$7:D$371)"

# Row 6: add F6 label + MIN formulas
$ws.Range("F6").Value = "min"
$ws.Range("G6").Formula = "=MIN(C$7:C$371)"
$ws.Range("H6").Formula = "=MIN(D$7:D$371)"

# Row 7: add F7 label + AVERAGE formulas
$ws.Range("F7").Value = "average"
$ws.Range("G7").Formula = "=AVERAGE(C$7:C$371)"
$ws.Range("H7").Formula = "=AVERAGE(D$7:D$371)"

# Move/resize the chart on this sheet to match new anchor position
$co = $ws.ChartObjects(1)
$co.Left = 350.625
$co.Top = 121.5
$co.Width = 806.562539
$co.Height = 256.125039

# Update the active selection to D26
$ws.Range("D26").Select()

Write-Output "done"
